$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "limon"
$ws.Range("B7").Value = 5
$ws.Range("A8").Value = "merengue"
$ws.Range("B8").Value = 2
$ws.Range("A9").Value = "vainilla"
$ws.Range("B9").Value = 1
$ws.Range("A10").Value = "leche"
$ws.Range("B10").Value = 8
